$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume(1h) data per latest scrape
# For Price column (D), force Text number format before assignment so
# Excel does not auto-convert numeric-looking strings (e.g. "26.112.72")
# into actual numbers, matching the source data which stores them as text.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.112.72"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.665.59"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.71"
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("E8").Value = "  -2.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06317"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("E10").Value = "  -0.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07533"
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.647.35"
$ws.Range("E12").Value = "  -1.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.407"
$ws.Range("E13").Value = "  -2.23%  "
$ws.Range("E14").Value = "  -4.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000007994"
$ws.Range("E15").Value = "  -1.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.27"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.189.02"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.727"
$ws.Range("E19").Value = "  -2.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "186.88"
$ws.Range("E20").Value = "  -1.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.24"
$ws.Range("E21").Value = "  -3.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.216"
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.59"
$ws.Range("E24").Value = "  +0.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1237"
$ws.Range("E25").Value = "  -1.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.424"
$ws.Range("E26").Value = "  -2.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.71"
$ws.Range("E27").Value = "  -2.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06265"
$ws.Range("E28").Value = "  -1.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.360"
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.276"
$ws.Range("E30").Value = "  -0.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.491"
$ws.Range("E31").Value = "  -1.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.401"
$ws.Range("E32").Value = "  -3.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.636"
$ws.Range("E33").Value = "  -2.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9985"
$ws.Range("E34").Value = "  -1.19%  "
$ws.Range("B35").Value = "MXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.762"
$ws.Range("E35").Value = "  +1.63%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.393"
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5982"
$ws.Range("E37").Value = "  -1.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.109.04"
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01608"
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.050"
$ws.Range("E40").Value = "  -1.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8622"
$ws.Range("E41").Value = "  -1.10%  "
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.60"
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.815.78"
$ws.Range("E44").Value = "  -0.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000108"
$ws.Range("E45").Value = "  -1.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.26"
$ws.Range("E46").Value = "  -3.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.002"
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.073"
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05240"
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4232"
$ws.Range("E50").Value = "  -0.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.872"
$ws.Range("E51").Value = "  -1.39%  "
